$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) stays text, matching the source data which uses
# localized "thousand.thousand.decimal"-style strings that Excel would
# otherwise auto-convert to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "28.447.53"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.876.54"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.022"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "317.43"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.019"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.5121"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3959"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.08462"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "42.02"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.253"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.883.25"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "20.47"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.240"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.022"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "91.18"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.06778"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.69"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.019"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.935"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.505.07"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.087.40"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "161.88"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.351"
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "127.03"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1056"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.045"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.773"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02434"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06481"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2177"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.867"
$ws.Range("E38").Value = "  -6.67%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.265"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.185"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6393"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "4.999"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.21"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.6039"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.05"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.723"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "1.991"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "1.206"
$ws.Range("E48").Value = "  -6.41%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "122.06"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06842"
$ws.Range("E51").Value = "  -1.01%  "
